$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # used only as a source of already-existing shared strings
$ws7 = $wb.Worksheets.Item(7)   # 保險 (insurance)
$ws8 = $wb.Worksheets.Item(8)   # 債權 (claim / debt)

$xlPasteValues = -4163

# ---------------------------------------------------------------------------
# Sheet "保險" (insurance) -- 7th worksheet
# ---------------------------------------------------------------------------

# Header row (row 1): relabel existing columns, then add the new header cells.
# Copying from E1 first keeps the bold/bordered header style (s="1") on the
# newly created cells; the Value assignment afterwards overwrites the copied
# text with the real header label.
$ws7.Range("B1").Value = "company"
$ws7.Range("C1").Value = "name"
$ws7.Range("D1").Value = "owner"
$ws7.Range("E1").Value = "property_category"

$ws7.Range("E1").Copy($ws7.Range("F1")); $ws7.Range("F1").Value = "category"
$ws7.Range("E1").Copy($ws7.Range("G1")); $ws7.Range("G1").Value = "date"
$ws7.Range("E1").Copy($ws7.Range("H1")); $ws7.Range("H1").Value = "legislator_name"
$ws7.Range("E1").Copy($ws7.Range("I1")); $ws7.Range("I1").Value = "legislator_id"
$ws7.Range("E1").Copy($ws7.Range("J1")); $ws7.Range("J1").Value = "source_file"
$ws7.Range("E1").Copy($ws7.Range("K1")); $ws7.Range("K1").Value = "index"

# Data rows 2-9: column E becomes "insurance"; columns F..K are brand new and
# hold the same boilerplate metadata used throughout the rest of the workbook.
# For values that already exist verbatim elsewhere in the workbook (normal,
# 2012-04-20, 蔣乃辛, tmp7091) we copy+paste-values from such a cell instead
# of typing the literal text, because typing "2012-04-20" directly would be
# auto-parsed by Excel into a date serial number.
$indexValues = @{2 = 84; 3 = 85; 4 = 86; 5 = 87; 6 = 88; 7 = 89; 8 = 90; 9 = 91}

foreach ($r in 2..9) {
    $ws7.Range("E$r").Value = "insurance"

    $ws1.Range("J2").Copy()
    $ws7.Range("F$r").PasteSpecial($xlPasteValues)

    $ws1.Range("K2").Copy()
    $ws7.Range("G$r").PasteSpecial($xlPasteValues)

    $ws1.Range("L2").Copy()
    $ws7.Range("H$r").PasteSpecial($xlPasteValues)

    $ws1.Range("M2").Copy()
    $ws7.Range("I$r").PasteSpecial($xlPasteValues)

    $ws1.Range("N2").Copy()
    $ws7.Range("J$r").PasteSpecial($xlPasteValues)

    $ws7.Range("K$r").Value = $indexValues[$r]
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "債權" (claim / debt) -- 8th worksheet
# ---------------------------------------------------------------------------

# Header row (row 1): relabel existing columns, then add the new header cells.
$ws8.Range("B1").Value = "species"
$ws8.Range("C1").Value = "owner"
$ws8.Range("D1").Value = "debtor"
$ws8.Range("E1").Value = "total"
$ws8.Range("F1").Value = "register_date"
$ws8.Range("G1").Value = "register_reason"

$ws8.Range("G1").Copy($ws8.Range("H1")); $ws8.Range("H1").Value = "property_category"
$ws8.Range("G1").Copy($ws8.Range("I1")); $ws8.Range("I1").Value = "category"
$ws8.Range("G1").Copy($ws8.Range("J1")); $ws8.Range("J1").Value = "date"
$ws8.Range("G1").Copy($ws8.Range("K1")); $ws8.Range("K1").Value = "legislator_name"
$ws8.Range("G1").Copy($ws8.Range("L1")); $ws8.Range("L1").Value = "legislator_id"
$ws8.Range("G1").Copy($ws8.Range("M1")); $ws8.Range("M1").Value = "source_file"
$ws8.Range("G1").Copy($ws8.Range("N1")); $ws8.Range("N1").Value = "index"

# Data row 2: columns H..N are brand new.
$ws8.Range("H2").Value = "claim"

$ws1.Range("J2").Copy()
$ws8.Range("I2").PasteSpecial($xlPasteValues)

$ws1.Range("K2").Copy()
$ws8.Range("J2").PasteSpecial($xlPasteValues)

$ws1.Range("L2").Copy()
$ws8.Range("K2").PasteSpecial($xlPasteValues)

$ws1.Range("M2").Copy()
$ws8.Range("L2").PasteSpecial($xlPasteValues)

$ws1.Range("N2").Copy()
$ws8.Range("M2").PasteSpecial($xlPasteValues)

$ws8.Range("N2").Value = 96
$excel.CutCopyMode = $false
